$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 102, shifting every existing row (102-143)
# down by one (new last row becomes 144).
$ws.Rows.Item(102).Insert()

# Populate the newly inserted row 102 with the new weekly record.
$ws.Cells.Item(102, 1).Value = 10
$ws.Cells.Item(102, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(102, 3).Value = "La Araucanía"
$ws.Cells.Item(102, 4).Value = 44636
$ws.Cells.Item(102, 5).Value = 9
$ws.Cells.Item(102, 6).Value = 100114007
$ws.Cells.Item(102, 7).Value = "Jengibre"
$ws.Cells.Item(102, 8).Value = "Sin especificar"
$ws.Cells.Item(102, 9).Value = "Primera"
$ws.Cells.Item(102, 10).Value = 40
$ws.Cells.Item(102, 11).Value = 25000
$ws.Cells.Item(102, 12).Value = 25000
$ws.Cells.Item(102, 13).Value = 25000
$ws.Cells.Item(102, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(102, 15).Value = "Perú"
$ws.Cells.Item(102, 16).Value = 1923
$ws.Cells.Item(102, 17).Value = 13
$ws.Cells.Item(102, 18).Value = "Hortaliza"
